# Regenerate s_vals data to filter save games.
# Updates the numeric TB/d2S/K/IP/sum values (columns B-E, G) for rows 2-6.
# Column A (dates) and F (Win) are left unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(3.182878228561681, 1.65323645889881,  0.1529057820181812, 0.4998867070740569, 5.488907176552729)
    3 = @(3.182878228561681, 86.29678392075563, 3.082599426703578,  6.48142807727062,   99.04368965329151)
    4 = @(3.182878228561681, 1.65323645889881,  0.1529057820181812, 0.4998867070740569, 5.488907176552729)
    5 = @(3.182878228561681, 1.65323645889881,  3.082599426703578,  0.4998867070740569, 8.418600821238126)
    6 = @(1.505614041169197, 1.65323645889881,  3.082599426703578,  0.4998867070740569, 6.741336633845642)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]  # B - TB
    $ws.Cells.Item($row, 3).Value = $vals[1]  # C - d2S
    $ws.Cells.Item($row, 4).Value = $vals[2]  # D - K
    $ws.Cells.Item($row, 5).Value = $vals[3]  # E - IP
    $ws.Cells.Item($row, 7).Value = $vals[4]  # G - sum
}
